# New point maps for the new experiments
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the point map values for BMW_X5_SUV_1/Localization (row 5)
# and BMW_X5_SUV_1/FakeLocalization (row 6)
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 1

# Update the active selection to match the new edit location
$ws.Range("B7").Select()
